$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price cells that would otherwise be
# auto-converted to numbers by Excel (losing trailing zeros).
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'

$ws.Range('D2').Value = '32.717.58'
$ws.Range('E2').Value = '  +9.29%  '
$ws.Range('D3').Value = '1.756.17'
$ws.Range('E3').Value = '  +5.77%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '227.05'
$ws.Range('E5').Value = '  +4.43%  '
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '31.69'
$ws.Range('E8').Value = '  +9.73%  '
$ws.Range('D9').Value = '45.19'
$ws.Range('E9').Value = '  +3.01%  '
$ws.Range('D10').Value = '0.277'
$ws.Range('E10').Value = '  +5.17%  '
$ws.Range('D11').Value = '0.0666'
$ws.Range('E11').Value = '  +8.33%  '
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = '2.009.80'
$ws.Range('E13').Value = '  +5.92%  '
$ws.Range('D14').Value = '1.757.57'
$ws.Range('E14').Value = '  +5.92%  '
$ws.Range('E15').Value = '  +3.79%  '
$ws.Range('D16').Value = '10.42'
$ws.Range('E16').Value = '  +3.33%  '
$ws.Range('D17').Value = '4.28'
$ws.Range('E17').Value = '  +8.24%  '
$ws.Range('D18').Value = '32.796.72'
$ws.Range('E18').Value = '  +9.46%  '
$ws.Range('D19').Value = '68.72'
$ws.Range('E19').Value = '  +5.67%  '
$ws.Range('D20').Value = '258.14'
$ws.Range('E20').Value = '  +6.23%  '
$ws.Range('E21').Value = '  +4.37%  '
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').Value = '10.43'
$ws.Range('E23').Value = '  +3.92%  '
$ws.Range('D24').Value = '4.35'
$ws.Range('E24').Value = '  +3.89%  '
$ws.Range('D25').Value = '2.17'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').Value = '159.80'
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('D27').Value = '16.52'
$ws.Range('E27').Value = '  +4.54%  '
$ws.Range('E28').Value = '  +4.03%  '
$ws.Range('D29').Value = '6.96'
$ws.Range('E29').Value = '  +3.38%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').Value = '3.88'
$ws.Range('E31').Value = '  +13.90%  '
$ws.Range('E32').Value = '  +3.22%  '
$ws.Range('E33').Value = '  +5.32%  '
$ws.Range('D34').Value = '3.49'
$ws.Range('E34').Value = '  +8.12%  '
$ws.Range('D35').Value = '1.552.11'
$ws.Range('E35').Value = '  +7.38%  '
$ws.Range('D36').Value = '1.78'
$ws.Range('E36').Value = '  +3.95%  '
$ws.Range('D37').Value = '1.04'
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('E38').Value = '  +9.84%  '
$ws.Range('D39').Value = '84.31'
$ws.Range('E39').Value = '  +6.02%  '
$ws.Range('D40').Value = '0.0185'
$ws.Range('E40').Value = '  +5.45%  '
$ws.Range('E41').Value = '  +3.56%  '
$ws.Range('D42').Value = '2.31'
$ws.Range('E42').Value = '  +0.85%  '
$ws.Range('D43').Value = '0.872'
$ws.Range('E43').Value = '  +2.56%  '
$ws.Range('D44').Value = '2.07'
$ws.Range('E44').Value = '  +6.18%  '
$ws.Range('D45').Value = '0.0515'
$ws.Range('E45').Value = '  +2.80%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = '1.05'
$ws.Range('E46').Value = '  +4.29%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '54.34'
$ws.Range('E47').Value = '  +7.52%  '
$ws.Range('E48').Value = '  +5.82%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = '5.67'
$ws.Range('E50').Value = '  +5.67%  '
$ws.Range('D51').Value = '95.49'
$ws.Range('E51').Value = '  +1.36%  '
